$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Weekly crime statistics updates (rows 14-30) ---

# Row 14
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 50
$ws.Range("I14").Value = 17
$ws.Range("J14").Value = 13
$ws.Range("K14").Value = 30.769230769230
$ws.Range("L14").Value = -5.555555555555
$ws.Range("M14").Value = -26.086956521739
$ws.Range("N14").Value = -80.459770114942

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 4
$ws.Range("I15").Value = 43
$ws.Range("J15").Value = 33
$ws.Range("K15").Value = 30.303030303030
$ws.Range("L15").Value = 16.216216216216
$ws.Range("M15").Value = -2.272727272727
$ws.Range("N15").Value = -46.913580246913

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("E16").Value = -37.5
$ws.Range("F16").Value = 49
$ws.Range("G16").Value = 57
$ws.Range("H16").Value = -14.035087719298
$ws.Range("I16").Value = 452
$ws.Range("J16").Value = 527
$ws.Range("K16").Value = -14.231499051233
$ws.Range("L16").Value = 16.494845360824
$ws.Range("M16").Value = -8.686868686868
$ws.Range("N16").Value = -77.067478437341

# Row 17
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 26
$ws.Range("E17").Value = -34.615384615384
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 84
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 725
$ws.Range("J17").Value = 760
$ws.Range("K17").Value = -4.605263157894
$ws.Range("L17").Value = 24.784853700516
$ws.Range("M17").Value = 42.436149312377
$ws.Range("N17").Value = -30.622009569378

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -28.125
$ws.Range("I18").Value = 258
$ws.Range("J18").Value = 271
$ws.Range("K18").Value = -4.797047970479
$ws.Range("L18").Value = -3.370786516853
$ws.Range("M18").Value = -19.626168224299
$ws.Range("N18").Value = -78.955954323001

# Row 19
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 104
$ws.Range("H19").Value = -31.730769230769
$ws.Range("I19").Value = 734
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = -18.444444444444
$ws.Range("L19").Value = 7.624633431085
$ws.Range("M19").Value = 55.838641188959
$ws.Range("N19").Value = 15.047021943573

# Row 20
$ws.Range("C20").Value = 8
$ws.Range("E20").Value = -11.111111111111
$ws.Range("F20").Value = 44
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = 41.935483870967
$ws.Range("I20").Value = 363
$ws.Range("J20").Value = 324
$ws.Range("K20").Value = 12.037037037037
$ws.Range("L20").Value = 2.253521126760
$ws.Range("M20").Value = 80.597014925373
$ws.Range("N20").Value = -79.537767756482

# Row 21
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 81
$ws.Range("E21").Value = -19.753086419753
$ws.Range("F21").Value = 264
$ws.Range("G21").Value = 314
$ws.Range("H21").Value = -15.923566878980
$ws.Range("I21").Value = 2592
$ws.Range("J21").Value = 2828
$ws.Range("K21").Value = -8.345120226308
$ws.Range("L21").Value = 11.340206185567
$ws.Range("M21").Value = 25.581395348837
$ws.Range("N21").Value = -62.005277044854

# Row 22
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 30
$ws.Range("J22").Value = 53
$ws.Range("K22").Value = -43.396226415094
$ws.Range("L22").Value = 11.111111111111
$ws.Range("M22").Value = -34.782608695652

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -55.555555555555
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 34
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 228
$ws.Range("J23").Value = 242
$ws.Range("K23").Value = -5.785123966942
$ws.Range("L23").Value = 19.371727748691
$ws.Range("M23").Value = 56.164383561643

# Row 24
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 54
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 158
$ws.Range("G24").Value = 231
$ws.Range("H24").Value = -31.601731601731
$ws.Range("I24").Value = 1597
$ws.Range("J24").Value = 1967
$ws.Range("K24").Value = -18.810371123538
$ws.Range("L24").Value = 7.181208053691
$ws.Range("M24").Value = 59.381237524950

# Row 25
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 39
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 107
$ws.Range("H25").Value = -19.626168224299
$ws.Range("I25").Value = 854
$ws.Range("J25").Value = 898
$ws.Range("K25").Value = -4.899777282850
$ws.Range("L25").Value = 27.653213751868
$ws.Range("M25").Value = -32.861635220125

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 65
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 30
$ws.Range("L26").Value = 0

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 85
$ws.Range("K27").Value = -5.882352941176
$ws.Range("L27").Value = -8.045977011494

# Row 28
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 48
$ws.Range("J28").Value = 67
$ws.Range("K28").Value = -28.358208955223
$ws.Range("L28").Value = -22.580645161290
$ws.Range("M28").Value = -34.246575342465
$ws.Range("N28").Value = -81.25

# Row 29
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -50
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 42
$ws.Range("J29").Value = 58
$ws.Range("K29").Value = -27.586206896551
$ws.Range("L29").Value = -14.285714285714
$ws.Range("M29").Value = -28.813559322033
$ws.Range("N29").Value = -81.739130434782

# Row 30
$ws.Range("L30").Value = -33.333333333333
